$wb = $excel.ActiveWorkbook

# --- Overview sheet: status text for both locale columns (rows 2 & 3) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value() = "Handed back: not in sync with en-US"
$overview.Range("F2").Value() = "Handed back: not in sync with en-US"
$overview.Range("E3").Value() = "Handed back: not in sync with en-US"
$overview.Range("F3").Value() = "Handed back: not in sync with en-US"

# --- zh-cn sheet: status text + refreshed handback datetime for row 3 ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value() = "Handed back: not in sync with en-US"
$zhcn.Range("C3").Value() = "Handed back: not in sync with en-US"
$zhcn.Range("K3").Value() = "2016-11-02 04:58:54"

# --- de-de sheet: status text + refreshed handback datetime for row 3 ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value() = "Handed back: not in sync with en-US"
$dede.Range("C3").Value() = "Handed back: not in sync with en-US"
$dede.Range("K3").Value() = "2016-11-02 04:59:12"

# --- widen the status columns to fit the longer text (report regeneration) ---
$overview.Columns("E").ColumnWidth = 32.6
$overview.Columns("F").ColumnWidth = 32.6
$zhcn.Columns("C").ColumnWidth = 32.6
$dede.Columns("C").ColumnWidth = 32.6
